$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.419420547413259
$ws.Range("B3").Value = 4.411282428249164
$ws.Range("B4").Value = 4.40770415605968
$ws.Range("B5").Value = 4.406607896379041
$ws.Range("B6").Value = 4.406447867465733
$ws.Range("B7").Value = 4.407687899504357
$ws.Range("B8").Value = 4.416323329831891
$ws.Range("B9").Value = 4.444267487759928
$ws.Range("B10").Value = 4.47115760396898
$ws.Range("B11").Value = 4.484690186219558
$ws.Range("B12").Value = 4.489994386205242
$ws.Range("B13").Value = 4.488844157273171
$ws.Range("B14").Value = 4.485123006728779
$ws.Range("B15").Value = 4.482866878356711
$ws.Range("B16").Value = 4.470298765143369
$ws.Range("B17").Value = 4.462916543587588
$ws.Range("B18").Value = 4.45879330924845
$ws.Range("B19").Value = 4.457418565337323
$ws.Range("B20").Value = 4.46368973364733
$ws.Range("B21").Value = 4.486211180619962
$ws.Range("B22").Value = 4.501974189432904
$ws.Range("B23").Value = 4.493468106653324
$ws.Range("B24").Value = 4.463339797287854
$ws.Range("B25").Value = 4.435564323666022
$ws.Range("D2").Value = 6.836331367955746
$ws.Range("D3").Value = 6.820337719359035
$ws.Range("D4").Value = 6.811338030866946
$ws.Range("D5").Value = 6.807875731986678
$ws.Range("D6").Value = 6.80731317312708
$ws.Range("D7").Value = 6.811290508110159
$ws.Range("D8").Value = 6.830645805779064
$ws.Range("D9").Value = 6.875159080456676
$ws.Range("D10").Value = 6.911919033541094
$ws.Range("D11").Value = 6.92953231875834
$ws.Range("D12").Value = 6.936330583338708
$ws.Range("D13").Value = 6.934860736219848
$ws.Range("D14").Value = 6.930089039243112
$ws.Range("D15").Value = 6.927182979438006
$ws.Range("D16").Value = 6.910785928282287
$ws.Range("D17").Value = 6.900955386835804
$ws.Range("D18").Value = 6.89538497198038
$ws.Range("D19").Value = 6.893513313163408
$ws.Range("D20").Value = 6.901993183165445
$ws.Range("D21").Value = 6.931487112761281
$ws.Range("D22").Value = 6.951512138919084
$ws.Range("D23").Value = 6.940755784592127
$ws.Range("D24").Value = 6.901523742261023
$ws.Range("D25").Value = 6.862408383266469
$ws.Range("E2").Value = 16.34318627713774
$ws.Range("E3").Value = 15.41977829608477
$ws.Range("E4").Value = 14.82996885745559
$ws.Range("E5").Value = 14.58414972551237
$ws.Range("E6").Value = 14.54301052625733
$ws.Range("E7").Value = 14.82667538861879
$ws.Range("E8").Value = 16.0296832317659
$ws.Range("E9").Value = 18.23741248690471
$ws.Range("E10").Value = 19.88179303399366
$ws.Range("E11").Value = 20.58903507246299
$ws.Range("E12").Value = 20.85102546124868
$ws.Range("E13").Value = 20.79485966584222
$ws.Range("E14").Value = 20.61070573188105
$ws.Range("E15").Value = 20.49714847529762
$ws.Range("E16").Value = 19.83475557826926
$ws.Range("E17").Value = 19.41797851724328
$ws.Range("E18").Value = 19.17441720988512
$ws.Range("E19").Value = 19.0912904385633
$ws.Range("E20").Value = 19.46274251035256
$ws.Range("E21").Value = 20.66495403058932
$ws.Range("E22").Value = 21.4167238387344
$ws.Range("E23").Value = 21.01858309254999
$ws.Range("E24").Value = 19.44251701736302
$ws.Range("E25").Value = 17.63300241715513
$ws.Range("F2").Value = 39.08880397238646
$ws.Range("F3").Value = 38.3838525370019
$ws.Range("F4").Value = 37.9553529045369
$ws.Range("F5").Value = 37.78203837576969
$ws.Range("F6").Value = 37.75334442105089
$ws.Range("F7").Value = 37.95300997319097
$ws.Range("F8").Value = 38.84496168865904
$ws.Range("F9").Value = 40.6193032943863
$ws.Range("F10").Value = 41.92608975985789
$ws.Range("F11").Value = 42.51877656181893
$ws.Range("F12").Value = 42.74275709660142
$ws.Range("F13").Value = 42.69454232697476
$ws.Range("F14").Value = 42.53721385978516
$ws.Range("F15").Value = 42.44078026881547
$ws.Range("F16").Value = 41.88730310190893
$ws.Range("F17").Value = 41.54716159417332
$ws.Range("F18").Value = 41.35136550099159
$ws.Range("F19").Value = 41.28505184758232
$ws.Range("F20").Value = 41.58338783949183
$ws.Range("F21").Value = 42.58343901303962
$ws.Range("F22").Value = 43.23427880488854
$ws.Range("F23").Value = 42.88722942174376
$ws.Range("F24").Value = 41.5670106992529
$ws.Range("F25").Value = 40.1378710247814
$ws.Range("G2").Value = 3.663932796302124
$ws.Range("G3").Value = 3.668712872757238
$ws.Range("G4").Value = 3.671790731665832
$ws.Range("G5").Value = 3.673081088138693
$ws.Range("G6").Value = 3.673297536628231
$ws.Range("G7").Value = 3.671807987444707
$ws.Range("G8").Value = 3.665551420691172
$ws.Range("G9").Value = 3.654407813508791
$ws.Range("G10").Value = 3.646895313483844
$ws.Range("G11").Value = 3.643621721146215
$ws.Range("G12").Value = 3.642402593413872
$ws.Range("G13").Value = 3.642664245030211
$ws.Range("G14").Value = 3.643521012740805
$ws.Range("G15").Value = 3.644048473454911
$ws.Range("G16").Value = 3.64711213039298
$ws.Range("G17").Value = 3.649028308842142
$ws.Range("G18").Value = 3.650143997596684
$ws.Range("G19").Value = 3.650524083929211
$ws.Range("G20").Value = 3.648822926927498
$ws.Range("G21").Value = 3.643268804075464
$ws.Range("G22").Value = 3.63975833159647
$ws.Range("G23").Value = 3.641621064975943
$ws.Range("G24").Value = 3.648915736307393
$ws.Range("G25").Value = 3.657303141046314
$ws.Range("I2").Value = 23.69124698385301
$ws.Range("I3").Value = 23.58631559592771
$ws.Range("I4").Value = 23.5242663693759
$ws.Range("I5").Value = 23.49958353029885
$ws.Range("I6").Value = 23.49552148156157
$ws.Range("I7").Value = 23.52393104223736
$ws.Range("I8").Value = 23.65457578452222
$ws.Range("I9").Value = 23.9294207382217
$ws.Range("I10").Value = 24.14237779743867
$ws.Range("I11").Value = 24.2415643997835
$ws.Range("I12").Value = 24.27944733455672
$ws.Range("I13").Value = 24.27127434997344
$ws.Range("I14").Value = 24.24467463110227
$ws.Range("I15").Value = 24.22842332946034
$ws.Range("I16").Value = 24.1359414545756
$ws.Range("I17").Value = 24.07979244531262
$ws.Range("I18").Value = 24.04771517433131
$ws.Range("I19").Value = 24.03689210066756
$ws.Range("I20").Value = 24.08574709002179
$ws.Range("I21").Value = 24.2524789208376
$ws.Range("I22").Value = 24.36332646942407
$ws.Range("I23").Value = 24.30399637043971
$ws.Range("I24").Value = 24.08305436046597
$ws.Range("I25").Value = 23.85310337205034
$ws.Range("K2").Value = 15.63710161176168
$ws.Range("K3").Value = 15.18622839472495
$ws.Range("K4").Value = 14.90901789620638
$ws.Range("K5").Value = 14.79616008649377
$ws.Range("K6").Value = 14.77743269300011
$ws.Range("K7").Value = 14.90749514761629
$ws.Range("K8").Value = 15.48183583021356
$ws.Range("K9").Value = 16.59681288681447
$ws.Range("K10").Value = 17.39873479435061
$ws.Range("K11").Value = 17.75786821657916
$ws.Range("K12").Value = 17.89290508107472
$ws.Range("K13").Value = 17.86386729391852
$ws.Range("K14").Value = 17.76899778395051
$ws.Range("K15").Value = 17.7107584977106
$ws.Range("K16").Value = 17.37513792209304
$ws.Range("K17").Value = 17.16768750111224
$ws.Range("K18").Value = 17.04784208053309
$ws.Range("K19").Value = 17.00717886474928
$ws.Range("K20").Value = 17.18982642225953
$ws.Range("K21").Value = 17.79689037547132
$ws.Range("K22").Value = 18.18799339927259
$ws.Range("K23").Value = 17.97981546751313
$ws.Range("K24").Value = 17.17981920962569
$ws.Range("K25").Value = 16.29751028732146
